# Update countries & provincias Spain
# Applies the latest COVID-19 data refresh to the "Pais" sheet:
#  - updates the "last updated" timestamp in A1
#  - swaps the country name in column A for rows whose ranking
#    flipped with their neighbour
#  - refreshes the numeric columns (B:H) wherever the source data changed

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Pais")

# --- A1: refreshed timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 13 de Septiembre de 2020 a las 14:43"

# --- Column A: country-name swaps (ranking changed) --------------------------
$ws.Cells.Item(44, 1).Value = "Paises Bajos"
$ws.Cells.Item(45, 1).Value = "Guatemala"
$ws.Cells.Item(48, 1).Value = "Bielorrusia"
$ws.Cells.Item(49, 1).Value = "Polonia"
$ws.Cells.Item(85, 1).Value = "Republica de Macedonia"
$ws.Cells.Item(86, 1).Value = "Madagascar"
$ws.Cells.Item(129, 1).Value = "Somalia"
$ws.Cells.Item(130, 1).Value = "Gambia"
$ws.Cells.Item(147, 1).Value = "Benin"
$ws.Cells.Item(148, 1).Value = "Botsuana"

# --- Columns B:H: refreshed numeric figures -----------------------------------
$ws.Cells.Item(4, 2).Value = 6678149
$ws.Cells.Item(4, 3).Value = 1548
$ws.Cells.Item(4, 4).Value = 3950599
$ws.Cells.Item(4, 5).Value = 2529402
$ws.Cells.Item(4, 7).Value = 20
$ws.Cells.Item(4, 8).Value = 198148

$ws.Cells.Item(5, 2).Value = 4764786
$ws.Cells.Item(5, 3).Value = 12998
$ws.Cells.Item(5, 4).Value = 3708257
$ws.Cells.Item(5, 5).Value = 977802
$ws.Cells.Item(5, 7).Value = 113
$ws.Cells.Item(5, 8).Value = 78727

$ws.Cells.Item(31, 2).Value = 121740
$ws.Cells.Item(31, 3).Value = 217
$ws.Cells.Item(31, 4).Value = 118682
$ws.Cells.Item(31, 5).Value = 2853

$ws.Cells.Item(38, 2).Value = 94764
$ws.Cells.Item(38, 3).Value = 553
$ws.Cells.Item(38, 4).Value = 84995
$ws.Cells.Item(38, 5).Value = 9209
$ws.Cells.Item(38, 7).Value = 2
$ws.Cells.Item(38, 8).Value = 560

$ws.Cells.Item(44, 2).Value = 82099
$ws.Cells.Item(44, 3).Value = 1087
$ws.Cells.Item(44, 4).Value = 0
$ws.Cells.Item(44, 5).Value = 0
$ws.Cells.Item(44, 7).Value = 1
$ws.Cells.Item(44, 8).Value = 6254

$ws.Cells.Item(45, 2).Value = 81658
$ws.Cells.Item(45, 4).Value = 70403
$ws.Cells.Item(45, 5).Value = 8306
$ws.Cells.Item(45, 8).Value = 2949

$ws.Cells.Item(48, 2).Value = 74173
$ws.Cells.Item(48, 3).Value = 198
$ws.Cells.Item(48, 4).Value = 72584
$ws.Cells.Item(48, 5).Value = 839
$ws.Cells.Item(48, 8).Value = 750

$ws.Cells.Item(49, 2).Value = 74152
$ws.Cells.Item(49, 3).Value = 502
$ws.Cells.Item(49, 4).Value = 60659
$ws.Cells.Item(49, 5).Value = 11305
$ws.Cells.Item(49, 7).Value = 6
$ws.Cells.Item(49, 8).Value = 2188

$ws.Cells.Item(61, 2).Value = 47042
$ws.Cells.Item(61, 3).Value = 321
$ws.Cells.Item(61, 5).Value = 3143
$ws.Cells.Item(61, 7).Value = 4
$ws.Cells.Item(61, 8).Value = 388

$ws.Cells.Item(75, 2).Value = 26928
$ws.Cells.Item(75, 3).Value = 77
$ws.Cells.Item(75, 4).Value = 18397
$ws.Cells.Item(75, 5).Value = 7746

$ws.Cells.Item(78, 2).Value = 23465
$ws.Cells.Item(78, 3).Value = 327
$ws.Cells.Item(78, 4).Value = 16038
$ws.Cells.Item(78, 5).Value = 6731
$ws.Cells.Item(78, 7).Value = 6
$ws.Cells.Item(78, 8).Value = 696

$ws.Cells.Item(82, 2).Value = 19890
$ws.Cells.Item(82, 3).Value = 333
$ws.Cells.Item(82, 4).Value = 16333
$ws.Cells.Item(82, 5).Value = 2926
$ws.Cells.Item(82, 7).Value = 1
$ws.Cells.Item(82, 8).Value = 631

$ws.Cells.Item(85, 2).Value = 15791
$ws.Cells.Item(85, 3).Value = 97
$ws.Cells.Item(85, 4).Value = 13184
$ws.Cells.Item(85, 5).Value = 1959
$ws.Cells.Item(85, 7).Value = 2
$ws.Cells.Item(85, 8).Value = 648

$ws.Cells.Item(86, 2).Value = 15757
$ws.Cells.Item(86, 3).Value = 20
$ws.Cells.Item(86, 4).Value = 14368
$ws.Cells.Item(86, 5).Value = 1178
$ws.Cells.Item(86, 7).Value = 1
$ws.Cells.Item(86, 8).Value = 211

$ws.Cells.Item(129, 2).Value = 3389
$ws.Cells.Item(129, 3).Value = 13
$ws.Cells.Item(129, 4).Value = 2803
$ws.Cells.Item(129, 5).Value = 488
$ws.Cells.Item(129, 8).Value = 98

$ws.Cells.Item(130, 4).Value = 1617
$ws.Cells.Item(130, 5).Value = 1657
$ws.Cells.Item(130, 8).Value = 102

$ws.Cells.Item(147, 2).Value = 2267
$ws.Cells.Item(147, 3).Value = 25
$ws.Cells.Item(147, 4).Value = 1942
$ws.Cells.Item(147, 5).Value = 285
$ws.Cells.Item(147, 8).Value = 40

$ws.Cells.Item(148, 2).Value = 2252
$ws.Cells.Item(148, 4).Value = 546
$ws.Cells.Item(148, 5).Value = 1696
$ws.Cells.Item(148, 8).Value = 10
